# Accommodate new design of w3.events: rename the "Network" node label to
# "Swarm" and shrink its box accordingly (on the two slides that show the
# Block/BP/Chain/Network diagram), and nudge the "In W3, we defeat..."
# paragraph on slide 2 to its new position.

$p = $ppt.ActivePresentation

# EMU per point, used to convert absolute EMU coordinates from the target
# OOXML into the points that the PowerPoint object model expects for
# Left/Top/Width/Height. Round to 4 decimals so the value converts back to
# the exact integer EMU amount (avoids off-by-one rounding).
$emuPerPt = 12700

function EmuToPt([double]$emu) {
    return [Math]::Round($emu / $emuPerPt, 4)
}

# --- Slide 2 ---------------------------------------------------------
$s2 = $p.Slides.Item(2)

# The diagram group ("组合 7") is the first shape; its "Network" label is
# group item 5 (Block, BP, Chain, Network in that order).
$diagram2 = $s2.Shapes.Item(1)
$networkLabel2 = $diagram2.GroupItems.Item(5)
$networkLabel2.Width = EmuToPt 1073646
$networkLabel2.TextFrame.TextRange.Text = "Swarm"

# The descriptive paragraph shape ("文本框 4") moves to a new position;
# its size is unchanged.
$paragraph2 = $s2.Shapes.Item(2)
$paragraph2.Left = EmuToPt 1416540
$paragraph2.Top = EmuToPt 640515

# --- Slide 3 ---------------------------------------------------------
$s3 = $p.Slides.Item(3)

$diagram3 = $s3.Shapes.Item(1)
$networkLabel3 = $diagram3.GroupItems.Item(5)
$networkLabel3.Width = EmuToPt 1073646
$networkLabel3.TextFrame.TextRange.Text = "Swarm"
